$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Template row (374) has the exact same cell layout/styles we need to replicate
# for each new row: inlineStr text cells A-D, blank E, date-formatted numeric F/G, blank H.
$templateRow = 374

$newRows = @(
    @(375, 'Fucntionality_test_MUTAG_with_SVC_Zero-GED_linear', 'MUTAG', 'SVC_Zero-GED_linear', 'SVC_Zero-GED_linear_trained_on_MUTAG.joblib', 45915.90239358796, 45915.90238939815),
    @(376, 'Fucntionality_test_MUTAG_with_SVC_Zero-GED_linear', 'MUTAG', 'SVC_Zero-GED_linear', 'SVC_Zero-GED_linear_trained_on_MUTAG.joblib', 45915.90305, 45915.90283905093),
    @(377, 'Fucntionality_test_MUTAG_with_SVC_Zero-GED_linear', 'MUTAG', 'SVC_Zero-GED_linear', 'SVC_Zero-GED_linear_trained_on_MUTAG.joblib', 45915.91021627315, 45915.90968631944),
    @(378, 'Fucntionality_test_MUTAG_with_SVC_Zero-GED_linear', 'MUTAG', 'SVC_Zero-GED_linear', 'SVC_Zero-GED_linear_trained_on_MUTAG.joblib', 45915.91075240741, 45915.91060165509),
    @(379, 'Fucntionality_test_MUTAG_with_SVC_Zero-GED_linear', 'MUTAG', 'SVC_Zero-GED_linear', 'SVC_Zero-GED_linear_trained_on_MUTAG.joblib', 45915.91097590278, 45915.91097273148),
    @(380, 'Fucntionality_test_MUTAG_with_SVC_Zero-GED_linear', 'MUTAG', 'SVC_Zero-GED_linear', 'SVC_Zero-GED_linear_trained_on_MUTAG.joblib', 45915.92173537037, 45915.92173221065),
    @(381, 'Fucntionality_test_MUTAG_with_SVC_Zero-GED_linear', 'MUTAG', 'SVC_Zero-GED_linear', 'SVC_Zero-GED_linear_trained_on_MUTAG.joblib', 45915.9218982176, 45915.92189467593),
    @(382, 'Fucntionality_test_MUTAG_with_SVC_Zero-GED_linear', 'MUTAG', 'SVC_Zero-GED_linear', 'SVC_Zero-GED_linear_trained_on_MUTAG.joblib', 45915.92199268519, 45915.92198901621),
    @(383, 'Fucntionality_test_MUTAG_with_SVC_Zero-GED_linear', 'MUTAG', 'SVC_Zero-GED_linear', 'SVC_Zero-GED_linear_trained_on_MUTAG.joblib', 45915.92215949074, 45915.92215592592),
    @(384, 'Fucntionality_test_MUTAG_with_SVC_Zero-GED_linear', 'MUTAG', 'SVC_Zero-GED_linear', 'SVC_Zero-GED_linear_trained_on_MUTAG.joblib', 45915.92252958333, 45915.92252607639),
    @(385, 'Fucntionality_test_MUTAG_with_SVC_Zero-GED_linear', 'MUTAG', 'SVC_Zero-GED_linear', 'SVC_Zero-GED_linear_trained_on_MUTAG.joblib', 45915.92266288195, 45915.92265982639),
    @(386, 'Fucntionality_test_MUTAG_with_SVC_Zero-GED_linear', 'MUTAG', 'SVC_Zero-GED_linear', 'SVC_Zero-GED_linear_trained_on_MUTAG.joblib', 45915.92331371528, 45915.92331065972),
    @(387, 'Fucntionality_test_MUTAG_with_SVC_Zero-GED_linear', 'MUTAG', 'SVC_Zero-GED_linear', 'SVC_Zero-GED_linear_trained_on_MUTAG.joblib', 45915.92477596065, 45915.92477296296),
    @(388, 'Fucntionality_test_MUTAG_with_SVC_Zero-GED_linear', 'MUTAG', 'SVC_Zero-GED_linear', 'SVC_Zero-GED_linear_trained_on_MUTAG.joblib', 45915.92644039352, 45915.92643734954),
    @(389, 'Fucntionality_test_MUTAG_with_SVC_Zero-GED_linear', 'MUTAG', 'SVC_Zero-GED_linear', 'SVC_Zero-GED_linear_trained_on_MUTAG.joblib', 45915.92711929399, 45915.92711626158),
    @(390, 'Fucntionality_test_MUTAG_with_SVC_Zero-GED_linear', 'MUTAG', 'SVC_Zero-GED_linear', 'SVC_Zero-GED_linear_trained_on_MUTAG.joblib', 45915.92757193287, 45915.92756888889),
    @(391, 'Fucntionality_test_MUTAG_with_SVC_Zero-GED_linear', 'MUTAG', 'SVC_Zero-GED_linear', 'SVC_Zero-GED_linear_trained_on_MUTAG.joblib', 45915.92806930556, 45915.92806621527),
    @(392, 'Fucntionality_test_MUTAG_with_SVC_Zero-GED_linear', 'MUTAG', 'SVC_Zero-GED_linear', 'SVC_Zero-GED_linear_trained_on_MUTAG.joblib', 45915.93155193287, 45915.93094765046),
    @(393, 'Fucntionality_test_MUTAG_with_SVC_Zero-GED_linear', 'MUTAG', 'SVC_Zero-GED_linear', 'SVC_Zero-GED_linear_trained_on_MUTAG.joblib', 45915.95589414352, 45915.9558908912),
    @(394, 'Fucntionality_test_MUTAG_with_SVC_Zero-GED_linear', 'MUTAG', 'SVC_Zero-GED_linear', 'SVC_Zero-GED_linear_trained_on_MUTAG.joblib', 45915.9568825, 45915.95687934028),
    @(395, 'Fucntionality_test_MUTAG_with_SVC_Zero-GED_linear', 'MUTAG', 'SVC_Zero-GED_linear', 'SVC_Zero-GED_linear_trained_on_MUTAG.joblib', 45915.95900917824, 45915.95900587963),
    @(396, 'Fucntionality_test_MUTAG_with_SVC_Zero-GED_linear', 'MUTAG', 'SVC_Zero-GED_linear', 'SVC_Zero-GED_linear_trained_on_MUTAG.joblib', 45915.96002582176, 45915.96002253472),
    @(397, 'Fucntionality_test_MUTAG_with_SVC_Zero-GED_linear', 'MUTAG', 'SVC_Zero-GED_linear', 'SVC_Zero-GED_linear_trained_on_MUTAG.joblib', 45915.96197258102, 45915.96196956019),
    @(398, 'Fucntionality_test_MUTAG_with_SVC_Zero-GED_linear', 'MUTAG', 'SVC_Zero-GED_linear', 'SVC_Zero-GED_linear_trained_on_MUTAG.joblib', 45915.96344083334, 45915.96343762732),
    @(399, 'Fucntionality_test_MUTAG_with_SVC_Simple-Prototype-GED_poly', 'MUTAG', 'SVC_Simple-Prototype-GED_poly', 'SVC_Simple-Prototype-GED_poly_trained_on_MUTAG.joblib', 45915.97400869213, 45915.97400861111),
    @(400, 'Fucntionality_test_MUTAG_with_SVC_Simple-Prototype-GED_poly', 'MUTAG', 'SVC_Simple-Prototype-GED_poly', 'SVC_Simple-Prototype-GED_poly_trained_on_MUTAG.joblib', 45915.97696756945, 45915.97696747685),
    @(401, 'Fucntionality_test_MUTAG_with_SVC_Simple-Prototype-GED_poly', 'MUTAG', 'SVC_Simple-Prototype-GED_poly', 'SVC_Simple-Prototype-GED_poly_trained_on_MUTAG.joblib', 45915.97728130787, 45915.97723453704),
    @(402, 'Fucntionality_test_MUTAG_with_SVC_Simple-Prototype-GED_poly', 'MUTAG', 'SVC_Simple-Prototype-GED_poly', 'SVC_Simple-Prototype-GED_poly_trained_on_MUTAG.joblib', 45915.97757503473, 45915.97757493056),
    @(403, 'Fucntionality_test_MUTAG_with_SVC_Simple-Prototype-GED_poly', 'MUTAG', 'SVC_Simple-Prototype-GED_poly', 'SVC_Simple-Prototype-GED_poly_trained_on_MUTAG.joblib', 45915.97782832831, 45915.97782822531)
)

foreach ($row in $newRows) {
    $r = $row[0]
    $srcRange = $ws.Range("A" + $templateRow + ":H" + $templateRow)
    $dstRange = $ws.Range("A" + $r + ":H" + $r)
    $srcRange.Copy($dstRange)

    $ws.Cells.Item($r, 1).Value2 = $row[1]
    $ws.Cells.Item($r, 2).Value2 = $row[2]
    $ws.Cells.Item($r, 3).Value2 = $row[3]
    $ws.Cells.Item($r, 4).Value2 = $row[4]
    $ws.Cells.Item($r, 6).Value2 = $row[5]
    $ws.Cells.Item($r, 7).Value2 = $row[6]
}

# Row 374 timestamps were refreshed with slightly updated precision values
$ws.Cells.Item(374, 6).Value2 = 45915.6777375926
$ws.Cells.Item(374, 7).Value2 = 45915.67773736111

Write-Host ("UsedRange: " + $ws.UsedRange.Address())
